$d = $word.ActiveDocument

# 1. Multi-choice feeding trials paragraph: add panel references
$r1 = $d.Content.Find.Execute(
    "g DM/day of CP and NDF, respectively (Figure 2). The target intake of naive hares fell between the nutritional rails of Diets B and C.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "g DM/day of CP and NDF, respectively (Figure 2, panel A). The target intake of naive hares fell between the nutritional rails of Diets B and C (Figure 2, panel B).",
    2)
Write-Host "Replace 1: $r1"

# 2. Single-choice feeding trials: results by treatment paragraph - renumber Figure 3 -> Figure 2 panels C/D
$r2 = $d.Content.Find.Execute(
    "not significant (Figure 3, panel A). This pattern of intake rate resulted in hares on diets B and C to have CP and NDF intake rates closest to the target intake found in the multi-choice trials (Figure 3, panel B).",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "not significant (Figure 2, panel C). This pattern of intake rate resulted in hares on diets B and C to have CP and NDF intake rates closest to the target intake found in the multi-choice trials (Figure 2, panel D).",
    2)
Write-Host "Replace 2: $r2"

# 3. Weight change paragraph - Figure 4 -> Figure 3 (two occurrences)
$r3 = $d.Content.Find.Execute(
    "hares could only maintain their weight on diet B (16.79 %/day; Figure 4). The Tukey test showed that weight change differed significantly between diet A and all other diets (Figure 4).",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "hares could only maintain their weight on diet B (16.79 %/day; Figure 3). The Tukey test showed that weight change differed significantly between diet A and all other diets (Figure 3).",
    2)
Write-Host "Replace 3: $r3"

# 4. DMD paragraph - Figure 5 -> Figure 4 (two occurrences)
$r4 = $d.Content.Find.Execute(
    "CP digestibility increased significantly as diet CP increased from A to B (Figure 5). Diet A and B produced similar NDF digestive rates, which were significantly higher than that of diet C. There was no significant difference in NDF digestibility of diet D and the other three diets (Figure 5).",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "CP digestibility increased significantly as diet CP increased from A to B (Figure 4). Diet A and B produced similar NDF digestive rates, which were significantly higher than that of diet C. There was no significant difference in NDF digestibility of diet D and the other three diets (Figure 4).",
    2)
Write-Host "Replace 4: $r4"

# 5. GAM paragraph - append new sentences about heat maps / GAM results
$r5 = $d.Content.Find.Execute(
    "We also assessed how protein and fibre intake affected total DMD.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "We also assessed how protein and fibre intake affected total DMD. In terms of crude macronutrients, the heat map shows that hares begin to maintain weight as protein intake increases, at approximately 10 g DM/kg^0.75/day, but only at mid ranges of fibre intake. After protein intake reaches about 14 g DM/kg^0.75/day, fibre intake becomes irrelevant (Figure 5). At highest fibre intake rates relative to protein intake we see the highest weight loss. The GAM for this pattern showed that the interacting effects of protein and fibre intake were significant toward weight change and that this relationship was non-linear (edf > 10); together this interaction explained 29% of the deviation (Table 2). The effect of digestible macronutrient intake on weight change shows a similar response to that of crude macronutrient, but with highest performance occurring from a more balanced intake of digestible protein and digestible fibre (Figure 5). This is further expressed with a higher edf in the GAM output, which means the effect is more non-linear. Yogether this interaction explained 52% of the deviation (Table 2). Lastly, the relationship between DMD and protein-fibre intake is the most non-linear of the three models, with the highest edf and the highest deviation explained (81%; Table 2). Figure 5 displays this, with higher digestibility occurring as the ratio of protein to fibre intake increases.",
    2)
Write-Host "Replace 5: $r5"
